$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.386.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.330.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.14"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.30"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.57"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.741.89"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.364.80"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.339.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.43"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "323.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.14"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.70"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +12.14%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.162"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.41"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0720"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.67"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.09"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.26"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.93"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.886"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.46"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.41"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.377"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.57"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.89"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.03"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0925"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.556"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.92"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.381"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.14"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.09%  "
